$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 9
$ws.Range("D2").Value = 44874
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 40
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 25000
$ws.Range("P2").Value = 25000
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("R2").Value = "Provincia de Quillota"
$ws.Range("S2").Value = 2500
$ws.Range("T2").Value = 10

# Row 3 <- original row 10
$ws.Range("D3").Value = 44483
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 35
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("Q3").Value = "$/bandeja 5 kilos"
$ws.Range("R3").Value = "Provincia de Quillota"
$ws.Range("S3").Value = 2000
$ws.Range("T3").Value = 5

# Row 4 <- original row 13
$ws.Range("D4").Value = 44496
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 55
$ws.Range("N4").Value = 28000
$ws.Range("O4").Value = 28000
$ws.Range("P4").Value = 28000
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("R4").Value = "Provincia de Quillota"
$ws.Range("S4").Value = 2800
$ws.Range("T4").Value = 10

# Row 5 <- original row 12
$ws.Range("D5").Value = 44921
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 55
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = "$/bandeja 7 kilos"
$ws.Range("R5").Value = "Provincia de Quillota"
$ws.Range("S5").Value = 2143
$ws.Range("T5").Value = 7

# Row 6 <- original row 2
$ws.Range("D6").Value = 44902
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 90
$ws.Range("N6").Value = 25000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 25000
$ws.Range("Q6").Value = "$/bandeja 10 kilos"
$ws.Range("R6").Value = "Provincia de Quillota"
$ws.Range("S6").Value = 2500
$ws.Range("T6").Value = 10

# Row 7 <- original row 15
$ws.Range("D7").Value = 44503
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 28000
$ws.Range("O7").Value = 28000
$ws.Range("P7").Value = 28000
$ws.Range("Q7").Value = "$/bandeja 10 kilos"
$ws.Range("R7").Value = "Provincia de Quillota"
$ws.Range("S7").Value = 2800
$ws.Range("T7").Value = 10

# Row 8 <- original row 19
$ws.Range("D8").Value = 44858
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 90
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("Q8").Value = "$/bandeja 5 kilos"
$ws.Range("R8").Value = "Provincia de Quillota"
$ws.Range("S8").Value = 4000
$ws.Range("T8").Value = 5

# Row 9 <- original row 6
$ws.Range("D9").Value = 44859
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("Q9").Value = "$/bandeja 5 kilos"
$ws.Range("R9").Value = "Provincia de Quillota"
$ws.Range("S9").Value = 4000
$ws.Range("T9").Value = 5

# Row 10 <- original row 7
$ws.Range("D10").Value = 44879
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 25
$ws.Range("N10").Value = 30000
$ws.Range("O10").Value = 30000
$ws.Range("P10").Value = 30000
$ws.Range("Q10").Value = "$/bandeja 10 kilos"
$ws.Range("R10").Value = "Provincia de Quillota"
$ws.Range("S10").Value = 3000
$ws.Range("T10").Value = 10

# Row 11 <- original row 5
$ws.Range("D11").Value = 44515
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 80
$ws.Range("N11").Value = 28000
$ws.Range("O11").Value = 28000
$ws.Range("P11").Value = 28000
$ws.Range("Q11").Value = "$/bandeja 10 kilos"
$ws.Range("R11").Value = "Provincia de Los Andes"
$ws.Range("S11").Value = 2800
$ws.Range("T11").Value = 10

# Row 12 <- original row 17
$ws.Range("D12").Value = 44488
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 12000
$ws.Range("P12").Value = 12000
$ws.Range("Q12").Value = "$/bandeja 5 kilos"
$ws.Range("R12").Value = "La Ligua"
$ws.Range("S12").Value = 2400
$ws.Range("T12").Value = 5

# Row 13 <- original row 14
$ws.Range("D13").Value = 44868
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 30
$ws.Range("N13").Value = 14000
$ws.Range("O13").Value = 14000
$ws.Range("P13").Value = 14000
$ws.Range("Q13").Value = "$/bandeja 5 kilos"
$ws.Range("R13").Value = "Provincia de Quillota"
$ws.Range("S13").Value = 2800
$ws.Range("T13").Value = 5

# Row 14 <- original row 8
$ws.Range("D14").Value = 44889
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 50
$ws.Range("N14").Value = 30000
$ws.Range("O14").Value = 30000
$ws.Range("P14").Value = 30000
$ws.Range("Q14").Value = "$/bandeja 10 kilos"
$ws.Range("R14").Value = "Provincia de Quillota"
$ws.Range("S14").Value = 3000
$ws.Range("T14").Value = 10

# Row 15 <- original row 3
$ws.Range("D15").Value = 44511
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 45
$ws.Range("N15").Value = 28000
$ws.Range("O15").Value = 28000
$ws.Range("P15").Value = 28000
$ws.Range("Q15").Value = "$/bandeja 10 kilos"
$ws.Range("R15").Value = "Provincia de Los Andes"
$ws.Range("S15").Value = 2800
$ws.Range("T15").Value = 10

# Row 16 <- original row 4
$ws.Range("D16").Value = 44511
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 45
$ws.Range("N16").Value = 3200
$ws.Range("O16").Value = 3200
$ws.Range("P16").Value = 3200
$ws.Range("Q16").Value = "$/bandeja 10 kilos"
$ws.Range("R16").Value = "Provincia de Quillota"
$ws.Range("S16").Value = 320
$ws.Range("T16").Value = 10

# Row 17 <- original row 11
$ws.Range("D17").Value = 44166
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("Q17").Value = "$/caja 18 kilos"
$ws.Range("R17").Value = "La Ligua"
$ws.Range("S17").Value = 667
$ws.Range("T17").Value = 18

# Row 18 <- original row 16
$ws.Range("D18").Value = 44519
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 30
$ws.Range("N18").Value = 28000
$ws.Range("O18").Value = 28000
$ws.Range("P18").Value = 28000
$ws.Range("Q18").Value = "$/bandeja 10 kilos"
$ws.Range("R18").Value = "Provincia de Quillota"
$ws.Range("S18").Value = 2800
$ws.Range("T18").Value = 10

# Row 19 <- original row 20
$ws.Range("D19").Value = 44901
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 40
$ws.Range("N19").Value = 25000
$ws.Range("O19").Value = 25000
$ws.Range("P19").Value = 25000
$ws.Range("Q19").Value = "$/bandeja 10 kilos"
$ws.Range("R19").Value = "Provincia de Quillota"
$ws.Range("S19").Value = 2500
$ws.Range("T19").Value = 10

# Row 20 <- original row 18
$ws.Range("D20").Value = 44466
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 80
$ws.Range("N20").Value = 11000
$ws.Range("O20").Value = 11000
$ws.Range("P20").Value = 11000
$ws.Range("Q20").Value = "$/bandeja 5 kilos"
$ws.Range("R20").Value = "La Ligua"
$ws.Range("S20").Value = 2200
$ws.Range("T20").Value = 5
